$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ILMN")

# Row 4 - Inventory
$ws.Range("B4").Value = 364000000.0
$ws.Range("C4").Value = 372000000.0
$ws.Range("D4").Value = 415000000.0
$ws.Range("E4").Value = 435000000.0
$ws.Range("F4").Value = 384000000.0

# Row 14 - Accounts Payable
$ws.Range("B14").Value = 178000000.0
$ws.Range("C14").Value = 192000000.0
$ws.Range("D14").Value = 156000000.0
$ws.Range("E14").Value = 135000000.0
$ws.Range("F14").Value = 130000000.0

# Row 22 - Long Term Tax Liability (Deferred)
$ws.Range("C22").Value = -20000000.0
$ws.Range("D22").Value = -19000000.0
$ws.Range("E22").Value = -13000000.0
$ws.Range("F22").Value = -91000000.0

# Row 36 - Net Debt
$ws.Range("G36").Value = -2228000000.0

# Row 37 - Total Debt
$ws.Range("G37").Value = 1186000000.0
